$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 24382.5
$ws.Range("J57").Value = 24382.5
$ws.Range("L57").Value = 73147.5
$ws.Range("N57").Value = -74145.5

$ws.Range("H62").Value = 3183.2104
$ws.Range("I62").Value = 2700
$ws.Range("K62").Value = 2700
$ws.Range("M62").Value = -2076

$ws.Range("H65").Value = 3183.2104
$ws.Range("I65").Value = 2700
$ws.Range("K65").Value = 13500
$ws.Range("M65").Value = -10380

$ws.Range("H86").Value = 11450.4
$ws.Range("I86").Value = 1737.5
$ws.Range("J86").Value = 17925.666
$ws.Range("K86").Value = 1737.5
$ws.Range("L86").Value = 17925.666
$ws.Range("M86").Value = -614.5
$ws.Range("N86").Value = -20171.666

$ws.Range("H88").Value = 1067
$ws.Range("I88").Value = 850.75
$ws.Range("J88").Value = 1499.5
$ws.Range("K88").Value = 850.75
$ws.Range("L88").Value = 1499.5
$ws.Range("M88").Value = -444.75
$ws.Range("N88").Value = -2311.5

$ws.Range("H89").Value = 11450.4
$ws.Range("I89").Value = 1737.5
$ws.Range("J89").Value = 17925.666
$ws.Range("K89").Value = 8687.5
$ws.Range("L89").Value = 89628.33
$ws.Range("M89").Value = -3071.5
$ws.Range("N89").Value = -100860.33

$ws.Range("H91").Value = 1067
$ws.Range("I91").Value = 850.75
$ws.Range("J91").Value = 1499.5
$ws.Range("K91").Value = 850.75
$ws.Range("L91").Value = 1499.5
$ws.Range("M91").Value = 553.25
$ws.Range("N91").Value = -4307.5

$ws.Range("H112").Value = 1051.7693
$ws.Range("J112").Value = 1061.84
$ws.Range("L112").Value = 3185.52
$ws.Range("N112").Value = -5401.52

$ws.Range("H129").Value = 385297.56
$ws.Range("J129").Value = 476959.53
$ws.Range("L129").Value = 1430878.59
$ws.Range("N129").Value = -1440878.59

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5931.202
$ws.Range("I32").Value = 4690.7427
$ws.Range("J32").Value = 12133.5
$ws.Range("K32").Value = 4690.7427
$ws.Range("L32").Value = 12133.5
$ws.Range("M32").Value = -4403.7427
$ws.Range("N32").Value = -12707.5

$ws.Range("H63").Value = 3126815
$ws.Range("I63").Value = 2016.7778
$ws.Range("K63").Value = 2016.7778
$ws.Range("M63").Value = -1330.7778

$ws.Range("H66").Value = 3126815
$ws.Range("I66").Value = 2016.7778
$ws.Range("K66").Value = 10083.889
$ws.Range("M66").Value = -6651.889000000001

$ws.Range("H132").Value = 18701
$ws.Range("I132").Value = 2426.65
$ws.Range("J132").Value = 48290.727
$ws.Range("K132").Value = 7279.950000000001
$ws.Range("L132").Value = 144872.181
$ws.Range("M132").Value = -4749.950000000001
$ws.Range("N132").Value = -149932.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5318.7666
$ws.Range("I134").Value = 5390.5654
$ws.Range("J134").Value = 5082.857
$ws.Range("K134").Value = 16171.6962
$ws.Range("L134").Value = 15248.571
$ws.Range("M134").Value = -13636.6962
$ws.Range("N134").Value = -20318.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2250

$ws.Range("H31").Value = 3304.5806
$ws.Range("I31").Value = 827.05554
$ws.Range("K31").Value = 827.05554
$ws.Range("M31").Value = -532.05554

$ws.Range("H34").Value = 3304.5806
$ws.Range("I34").Value = 827.05554
$ws.Range("K34").Value = 827.05554
$ws.Range("M34").Value = -625.05554

$ws.Range("H99").Value = 3261.56
$ws.Range("I99").Value = 2761.1765
$ws.Range("J99").Value = 4324.875
$ws.Range("K99").Value = 2761.1765
$ws.Range("L99").Value = 4324.875
$ws.Range("M99").Value = -1263.1765
$ws.Range("N99").Value = -7320.875

$ws.Range("H113").Value = 2250

$ws.Range("H126").Value = 3261.56
$ws.Range("I126").Value = 2761.1765
$ws.Range("J126").Value = 4324.875
$ws.Range("K126").Value = 8283.529500000001
$ws.Range("L126").Value = 12974.625
$ws.Range("M126").Value = -5813.529500000001
$ws.Range("N126").Value = -17914.625

$ws.Range("H132").Value = 3138.5264
$ws.Range("I132").Value = 1644.2858
$ws.Range("J132").Value = 7322.4
$ws.Range("K132").Value = 4932.857400000001
$ws.Range("L132").Value = 21967.2
$ws.Range("M132").Value = -2402.857400000001
$ws.Range("N132").Value = -27027.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2427.8
$ws.Range("I109").Value = 606.1667
$ws.Range("J109").Value = 3642.2222
$ws.Range("K109").Value = 1818.5001
$ws.Range("L109").Value = 10926.6666
$ws.Range("M109").Value = -778.5001
$ws.Range("N109").Value = -13006.6666

$ws.Range("H131").Value = 691.8
$ws.Range("J131").Value = 716.95654
$ws.Range("L131").Value = 2150.86962
$ws.Range("N131").Value = -12230.86962

$ws.Range("H134").Value = 2575.4375
$ws.Range("I134").Value = 1819.0834
$ws.Range("J134").Value = 4844.5
$ws.Range("K134").Value = 5457.2502
$ws.Range("L134").Value = 14533.5
$ws.Range("M134").Value = -387.2502000000004
$ws.Range("N134").Value = -24673.5

$ws.Range("H139").Value = 2227.625
$ws.Range("I139").Value = 1395.375
$ws.Range("K139").Value = 4186.125
$ws.Range("M139").Value = 953.875

$ws.Range("H140").Value = 1815.7646
$ws.Range("I140").Value = 841.0952
$ws.Range("K140").Value = 2523.2856
$ws.Range("M140").Value = 2656.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 3666.6667
$ws.Range("I49").Value = 2000
$ws.Range("K49").Value = 2000
$ws.Range("M49").Value = -1816

$ws.Range("H113").Value = 5114.276
$ws.Range("I113").Value = 6166.3
$ws.Range("J113").Value = 2776.4443
$ws.Range("K113").Value = 6166.3
$ws.Range("L113").Value = 2776.4443
$ws.Range("M113").Value = -3996.3
$ws.Range("N113").Value = -7116.4443

$ws.Range("H132").Value = 42571.5
$ws.Range("I132").Value = 8445.333000000001
$ws.Range("J132").Value = 103998.6
$ws.Range("K132").Value = 25335.999
$ws.Range("L132").Value = 311995.8
$ws.Range("M132").Value = -22805.999
$ws.Range("N132").Value = -317055.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 29995
$ws.Range("J116").Value = 29995
$ws.Range("L116").Value = 29995
$ws.Range("N116").Value = -39173

$ws.Range("H132").Value = 3457.1
$ws.Range("I132").Value = 2318.6
$ws.Range("K132").Value = 6955.799999999999
$ws.Range("M132").Value = -4425.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 12500
$ws.Range("J55").Value = 12500
$ws.Range("L55").Value = 12500
$ws.Range("N55").Value = -13054

$ws.Range("H81").Value = 2205.3333
$ws.Range("J81").Value = 2653.7144
$ws.Range("L81").Value = 5307.4288
$ws.Range("N81").Value = -7429.4288

$ws.Range("H84").Value = 2205.3333
$ws.Range("J84").Value = 2653.7144
$ws.Range("L84").Value = 26537.144
$ws.Range("N84").Value = -37145.144

$ws.Range("H132").Value = 1328.122
$ws.Range("I132").Value = 1111.2963
$ws.Range("J132").Value = 1746.2858
$ws.Range("K132").Value = 3333.8889
$ws.Range("L132").Value = 5238.857400000001
$ws.Range("M132").Value = -803.8888999999999
$ws.Range("N132").Value = -10298.8574
